$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: V1.24 feedback entry ---
$ws.Range("A4").Value = "V1.24"
$ws.Range("B4").Value = "- Game generally works well`n- Viking chess didn't implement all game rules`n- Adding sound for piece movements"
$ws.Range("C4").Value = "- Viking chess has been correctly implemented as of V1.25`n- Sounds will added soon"

# --- Row 5: V1.25 feedback entry ---
$ws.Range("A5").Value = "V1.25"
$ws.Range("B5").Value = "-Unintuative UI`n- Viking chess AI doesn't work`n- Normal chess AI sometimes freezes`n- Sometimes players appear twice in lobbies"
$ws.Range("C5").Value = "- UI will be overhalled`n- Viking chess AI has already been fixed as of V1.27`n- Normal chess AI will be looked at`n- Duplicate player bug will be looked at"

# Copy formatting (styles + row height) from row 3 down to the two new rows,
# matching the same look used by the existing feedback rows. PasteSpecial
# with formats-only leaves the values already written above untouched.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$ws.Range("A3:C3").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item(4).RowHeight = 90
$ws.Rows.Item(5).RowHeight = 90

$ws.Range("C6").Select()
$ws.Application.ActiveWindow.ScrollRow = 3
